$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Semestre ideal:" value changes from EF-6 to EF-4 (cells B9 and C9)
[void]$ws.Cells.Replace("EF-6", "EF-4")

# Remove the last requisito row (LOM3257 - Mecanica Classica), which drops
# the sheet's used range from A1:C24 down to A1:C23
$ws.Rows(24).Delete()

Write-Output "done"
